$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22 previously carried a batch of blank "placeholder" cells
# (I22 and K22:R22) left over from a formatting/fill operation. The
# review for that row only ever used columns A-H and J, so clear the
# stray blanks out of the row.
$ws.Range("I22").ClearContents()
$ws.Range("K22:R22").ClearContents()

# Add the new reviewed/translated product row right below it.
$ws.Range("A23").Value = "6VA18358"
$ws.Range("B23").Value = "W-7 BROW MASTER STENCIAL KIT"
$ws.Range("C23").Value = "Consumo"
$ws.Range("D23").Value = "No Tiene PT - TRADUZIDO"
$ws.Range("E23").Value = "Tiene ES"
$ws.Range("F23").Value = "No Tiene IT - TRADOTTO"
$ws.Range("G23").Value = "'4"
$ws.Range("H23").Value = "UND"
$ws.Range("J23").Value = "Revisado y Traducido"
